$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 11 and 12: coin identities swap (TRON <-> Dogecoin) plus new price/volume values
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D11" "0.0875"
$ws.Range("E11").Value = "  +1.22%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D12" "0.138"
$ws.Range("E12").Value = "  +1.67%  "

# Remaining price/volume updates
Set-TextValue "D2" "52.372.59"
$ws.Range("E2").Value = "  +0.10%  "
Set-TextValue "D3" "2.938.57"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "357.84"
$ws.Range("E5").Value = "  +1.31%  "
Set-TextValue "D6" "110.34"
$ws.Range("E6").Value = "  -1.74%  "
Set-TextValue "D7" "0.572"
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("E9").Value = "  +0.56%  "
Set-TextValue "D10" "39.26"
$ws.Range("E10").Value = "  -1.78%  "
Set-TextValue "D13" "19.59"
$ws.Range("E13").Value = "  -1.48%  "
Set-TextValue "D14" "7.87"
$ws.Range("E14").Value = "  +0.92%  "
Set-TextValue "D15" "3.398.43"
$ws.Range("E15").Value = "  +0.94%  "
Set-TextValue "D16" "2.924.44"
$ws.Range("E16").Value = "  +0.87%  "
Set-TextValue "D17" "0.986"
$ws.Range("E17").Value = "  -1.69%  "
Set-TextValue "D18" "52.352.11"
$ws.Range("E18").Value = "  +0.07%  "
Set-TextValue "D19" "3.55"
$ws.Range("E19").Value = "  +7.39%  "
Set-TextValue "D20" "7.61"
$ws.Range("E20").Value = "  -0.54%  "
Set-TextValue "D21" "14.02"
$ws.Range("E21").Value = "  -1.25%  "
Set-TextValue "D22" "0.0₃0988"
$ws.Range("E22").Value = "  +0.85%  "
Set-TextValue "D23" "70.70"
$ws.Range("E23").Value = "  -0.11%  "
Set-TextValue "D24" "271.27"
$ws.Range("E24").Value = "  +0.78%  "
Set-TextValue "D25" "2.82"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("E26").Value = "  +4.77%  "
Set-TextValue "D27" "7.97"
$ws.Range("E27").Value = "  +21.25%  "
Set-TextValue "D28" "27.10"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +5.92%  "
$ws.Range("E31").Value = "  -1.50%  "
Set-TextValue "D32" "37.80"
$ws.Range("E32").Value = "  -0.61%  "
Set-TextValue "D33" "2.29"
$ws.Range("E33").Value = "  +1.12%  "
Set-TextValue "D34" "6.19"
$ws.Range("E34").Value = "  -2.00%  "
Set-TextValue "D35" "52.34"
$ws.Range("E35").Value = "  -1.83%  "
Set-TextValue "D36" "0.0446"
$ws.Range("E36").Value = "  -1.44%  "
Set-TextValue "D37" "0.999"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -4.94%  "
Set-TextValue "D39" "18.31"
$ws.Range("E39").Value = "  -4.11%  "
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("E42").Value = "  +2.98%  "
Set-TextValue "D43" "23.03"
$ws.Range("E43").Value = "  -0.82%  "
Set-TextValue "D44" "120.08"
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("E45").Value = "  -0.89%  "
Set-TextValue "D46" "3.46"
$ws.Range("E46").Value = "  -1.86%  "
Set-TextValue "D47" "2.46"
$ws.Range("E47").Value = "  -5.78%  "
Set-TextValue "D48" "2.133.48"
$ws.Range("E48").Value = "  -2.62%  "
Set-TextValue "D49" "0.250"
$ws.Range("E49").Value = "  -4.54%  "
Set-TextValue "D50" "0.0354"
$ws.Range("E50").Value = "  +2.29%  "
Set-TextValue "D51" "0.926"
$ws.Range("E51").Value = "  -3.89%  "
